{"js": "// Table 2 (Greenspace/Bluespace) modification:\n// Swap the \"MVPA_Quant_PA2\" block (rows 4-8, 0-indexed) with the\n// \"MVPA_Quant_PA1\" block (rows 9-13, 0-indexed), and swap the\n// \"MVPA min/week - Activity count\" label (row 20) with the\n// \"MVPA min/week - Machine learning\" label (row 21).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// --- Read the two 5-row x 2-col blocks that need to swap ---\nconst blockACells = []; // rows 4-8 (PA2 block)\nconst blockBCells = []; // rows 9-13 (PA1 block)\nfor (let i = 0; i < 5; i++) {\n  const rowACells = [];\n  const rowBCells = [];\n  for (let c = 0; c < 2; c++) {\n    const cellA = table.getCell(4 + i, c);\n    cellA.load(\"value\");\n    rowACells.push(cellA);\n\n    const cellB = table.getCell(9 + i, c);\n    cellB.load(\"value\");\n    rowBCells.push(cellB);\n  }\n  blockACells.push(rowACells);\n  blockBCells.push(rowBCells);\n}\n\n// --- Read the two row-label cells that need to swap ---\nconst activityLabelCell = table.getCell(20, 0);\nactivityLabelCell.load(\"value\");\nconst machineLabelCell = table.getCell(21, 0);\nmachineLabelCell.load(\"value\");\n\nawait context.sync();\n\n// Capture the plain-text values before we start overwriting anything.\nconst blockAValues = blockACells.map((row) => row.map((cell) => cell.value));\nconst blockBValues = blockBCells.map((row) => row.map((cell) => cell.value));\nconst activityLabelValue = activityLabelCell.value;\nconst machineLabelValue = machineLabelCell.value;\n\n// --- Write the swapped values back ---\nfor (let i = 0; i < 5; i++) {\n  for (let c = 0; c < 2; c++) {\n    table.getCell(4 + i, c).value = blockBValues[i][c];\n    table.getCell(9 + i, c).value = blockAValues[i][c];\n  }\n}\n\ntable.getCell(20, 0).value = machineLabelValue;\ntable.getCell(21, 0).value = activityLabelValue;\n\nawait context.sync();\n", "ps1": "# Table 2 (Greenspace/Bluespace) modification:\n# Swap the \"MVPA_Quant_PA2\" block (rows 5-9) with the \"MVPA_Quant_PA1\"\n# block (rows 10-14), and swap the \"MVPA min/week - Activity count\" label\n# (row 21) with the \"MVPA min/week - Machine learning\" label (row 22).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Swap the PA2 block (rows 5-9) with the PA1 block (rows 10-14) ---\n# Column 1 (labels) and column 2 (values) for each of the 5 paired rows.\nfor ($i = 0; $i -lt 5; $i++) {\n    $rowA = 5 + $i\n    $rowB = 10 + $i\n\n    for ($col = 1; $col -le 2; $col++) {\n        $textA = $t.Cell($rowA, $col).Range.Text\n        $textB = $t.Cell($rowB, $col).Range.Text\n        $t.Cell($rowA, $col).Range.Text = $textB\n        $t.Cell($rowB, $col).Range.Text = $textA\n    }\n}\n\n# --- Swap the two row labels further down the table ---\n$labelActivity = $t.Cell(21, 1).Range.Text\n$labelMachine = $t.Cell(22, 1).Range.Text\n$t.Cell(21, 1).Range.Text = $labelMachine\n$t.Cell(22, 1).Range.Text = $labelActivity\n\nWrite-Output \"done\"\n"}
